$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header row (bold, bordered, centered - same as H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the I0/IF data values for rows 2-69.
$ijValues = @{
    2 = @(6, 7)
    3 = @(8, 8)
    4 = @(8, 9)
    5 = @(6, 6)
    6 = @(6, 7)
    7 = @(7, 8)
    8 = @(6, 7)
    9 = @(9, 9)
    10 = @(5, 6)
    11 = @(9, 9)
    12 = @(8, 8)
    13 = @(10, 11)
    14 = @(8, 8)
    15 = @(5, 6)
    16 = @(7, 7)
    17 = @(7, 8)
    18 = @(7, 7)
    19 = @(6, 7)
    20 = @(6, 6)
    21 = @(7, 7)
    22 = @(9, 9)
    23 = @(8, 8)
    24 = @(8, 8)
    25 = @(7, 7)
    26 = @(9, 9)
    27 = @(8, 8)
    28 = @(8, 8)
    29 = @(9, 9)
    30 = @(9, 9)
    31 = @(9, 9)
    32 = @(7, 7)
    33 = @(9, 9)
    34 = @(7, 7)
    35 = @(10, 10)
    36 = @(6, 6)
    37 = @(6, 7)
    38 = @(4, 5)
    39 = @(7, 9)
    40 = @(7, 8)
    41 = @(6, 7)
    42 = @(8, 9)
    43 = @(7, 7)
    44 = @(9, 9)
    45 = @(6, 7)
    46 = @(7, 8)
    47 = @(8, 8)
    48 = @(6, 6)
    49 = @(8, 9)
    50 = @(6, 7)
    51 = @(7, 7)
    52 = @(6, 6)
    53 = @(8, 9)
    54 = @(8, 9)
    55 = @(4, 7)
    56 = @(7, 8)
    57 = @(6, 6)
    58 = @(8, 8)
    59 = @(7, 7)
    60 = @(10, 10)
    61 = @(6, 7)
    62 = @(5, 6)
    63 = @(8, 8)
    64 = @(5, 6)
    65 = @(4, 4)
    66 = @(9, 9)
    67 = @(8, 8)
    68 = @(7, 7)
    69 = @(4, 4)
}

foreach ($r in $ijValues.Keys) {
    $pair = $ijValues[$r]
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
}

Write-Output "done"
